$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 13: new development log entry (2022-03-11, 10:30, 1 hour, OO Design + Implementation) ---

# C13: date, reuse the same date format as the rows above it (C10:C12)
$ws.Range("C12").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C13").Value2 = 44631

# D13: time value - clear any stale formatting first so Excel derives a fresh
# "time" number format for the typed value, same as it would for a blank cell.
$ws.Range("D13").ClearFormats()
$ws.Range("D13").Value2 = 0.4375
$ws.Range("D13").NumberFormat = "h:mm"

# E13: duration in hours
$ws.Range("E13").Value2 = 1

# F13 / G13: activity description (existing shared string "OO Design + Implementation")
$ws.Range("F13").Value = "OO Design + Implementation"
$ws.Range("G13").Value = "OO Design + Implementation"

[void]$ws.Range("G18").Select()
